$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.112.01"
$ws.Range("E2").Value = "  +0.23%  "

$ws.Range("D3").Value = "2.741.90"
$ws.Range("E3").Value = "  -0.43%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'570.19"
$ws.Range("E5").Value = "  -1.23%  "

$ws.Range("D6").Value = "'160.27"
$ws.Range("E6").Value = "  +1.62%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -1.49%  "

$ws.Range("D9").Value = "'0.110"
$ws.Range("E9").Value = "  -0.84%  "

$ws.Range("E10").Value = "  +4.79%  "

$ws.Range("D11").Value = "'5.81"
$ws.Range("E11").Value = "  +0.84%  "

$ws.Range("E12").Value = "  -0.19%  "

$ws.Range("D13").Value = "3.222.21"
$ws.Range("E13").Value = "  -0.53%  "

$ws.Range("D14").Value = "'26.87"
$ws.Range("E14").Value = "  +0.21%  "

$ws.Range("D15").Value = "63.926.15"
$ws.Range("E15").Value = "  +0.43%  "

$ws.Range("E16").Value = "  -0.81%  "

$ws.Range("D17").Value = "2.744.88"
$ws.Range("E17").Value = "  -0.44%  "

$ws.Range("D18").Value = "'12.14"
$ws.Range("E18").Value = "  +0.65%  "

$ws.Range("D19").Value = "'4.81"
$ws.Range("E19").Value = "  -1.10%  "

$ws.Range("D20").Value = "'354.50"
$ws.Range("E20").Value = "  -1.17%  "

$ws.Range("D21").Value = "'6.62"
$ws.Range("E21").Value = "  -2.39%  "

$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("E23").Value = "  -4.47%  "

$ws.Range("D24").Value = "'64.32"
$ws.Range("E24").Value = "  -2.17%  "

$ws.Range("E25").Value = "  +1.56%  "

$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("D27").Value = "'8.45"
$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("D28").Value = "0.0₃0922"
$ws.Range("E28").Value = "  -0.30%  "

$ws.Range("E29").Value = "  +1.15%  "

$ws.Range("E30").Value = "  +3.91%  "

$ws.Range("D31").Value = "'1.35"
$ws.Range("E31").Value = "  +9.50%  "

$ws.Range("D32").Value = "'164.27"
$ws.Range("E32").Value = "  -2.70%  "

$ws.Range("D33").Value = "'4.95"
$ws.Range("E33").Value = "  +0.48%  "

$ws.Range("D34").Value = "'20.09"
$ws.Range("E34").Value = "  -0.84%  "

$ws.Range("E35").Value = "  +2.23%  "

$ws.Range("D36").Value = "'0.998"
$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("D37").Value = "'1.82"
$ws.Range("E37").Value = "  +1.80%  "

$ws.Range("D38").Value = "'0.996"
$ws.Range("E38").Value = "  +0.45%  "

$ws.Range("D39").Value = "'350.15"
$ws.Range("E39").Value = "  +6.68%  "

$ws.Range("D40").Value = "'6.40"
$ws.Range("E40").Value = "  +4.51%  "

$ws.Range("E41").Value = "  -0.67%  "

$ws.Range("D42").Value = "'38.70"
$ws.Range("E42").Value = "  -1.47%  "

$ws.Range("E43").Value = "  +1.83%  "

$ws.Range("D44").Value = "'21.20"
$ws.Range("E44").Value = "  -1.40%  "

$ws.Range("D45").Value = "'0.0584"
$ws.Range("E45").Value = "  -1.04%  "

$ws.Range("D46").Value = "'0.627"
$ws.Range("E46").Value = "  -0.97%  "

$ws.Range("D47").Value = "'134.42"
$ws.Range("E47").Value = "  -0.84%  "

$ws.Range("D48").Value = "'0.101"
$ws.Range("E48").Value = "  -0.36%  "

$ws.Range("D49").Value = "'0.0251"
$ws.Range("E49").Value = "  -1.87%  "

$ws.Range("E50").Value = "  -0.12%  "

$ws.Range("D51").Value = "'11.06"
$ws.Range("E51").Value = "  +0.15%  "
